$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 4: Title "IETF 122 Meeting Tips" -> "IETF 123 Meeting Tips"
# ---------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$titleShape = $s4.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "IETF 123 Meeting Tips"

# ---------------------------------------------------------------------
# Slide 8: "TextBox 4" (the Outcome box)
# ---------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$tb = $s8.Shapes.Item(7)
$tr = $tb.TextFrame.TextRange

# --- paragraph 2: "dDesign team created, with authors ..." ---
#   remove the leading misspelled run "dDesign" entirely and let what
#   was " team " (unchanged formatting) become "Design team "
$para2 = $tr.Paragraphs(2, 1)
$para2.Characters(1, 7).Delete()              # removes the "dDesign" run completely
$firstRun = $para2.Runs(1, 1)
$firstRun.Text = "Design" + $firstRun.Text    # " team " -> "Design team "

# --- paragraph 3: "Nacho Dominguez ... Benoit Claise" ---
#   add a trailing space to the last run, then append "under Reshad supervision."
$para3 = $tr.Paragraphs(3, 1)
$lastRun = $para3.Runs(5, 1)
$lastRun.Text = $lastRun.Text + " "
$para3.InsertAfter("under") | Out-Null
$para3.InsertAfter(" Reshad supervision.") | Out-Null
